$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update title text (row 2) to reflect new reporting month
$ws.Range("A2").Value = "Commercial Sector by State, November 2016"

# Update West North Central row (row 22)
$ws.Range("B22").Value = 10
$ws.Range("C22").Value = 3.04
$ws.Range("D22").Value = 8.5

# Update Missouri row (row 26)
$ws.Range("B26").Value = 10
$ws.Range("C26").Value = 3.04
$ws.Range("D26").Value = 8.5

# Update U.S. Total row (row 66)
$ws.Range("B66").Value = 10
$ws.Range("C66").Value = 3.04
$ws.Range("D66").Value = 8.5
